# The sheet gains one new weekly price record for Albahaca (Vega Modelo de
# Temuco). It is inserted as the new row 168, pushing the former rows
# 168-212 down to become rows 169-213 (dimension grows from R212 to R213).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 168; Excel shifts rows 168:212 down to 169:213.
$ws.Rows.Item(168).Insert()

# Fill in the new record in row 168.
$ws.Cells.Item(168, 1).Value = 10
$ws.Cells.Item(168, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(168, 3).Value = "La Araucanía"
$ws.Cells.Item(168, 4).Value = 44627
$ws.Cells.Item(168, 5).Value = 9
$ws.Cells.Item(168, 6).Value = 100112052
$ws.Cells.Item(168, 7).Value = "Albahaca"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 15
$ws.Cells.Item(168, 11).Value = 5000
$ws.Cells.Item(168, 12).Value = 5000
$ws.Cells.Item(168, 13).Value = 5000
$ws.Cells.Item(168, 14).Value = "$/paquete"
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 5000
$ws.Cells.Item(168, 17).Value = 1
$ws.Cells.Item(168, 18).Value = "Hortaliza"
